$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking strings in the price/volume columns stay text,
# matching the workbook's original inline-string (t="inlineStr") storage.
$ws.Columns("B:E").NumberFormat = "@"

$ws.Range("D2").Value = "30.643.33"
$ws.Range("D3").Value = "1.937.93"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "246.31"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4817"
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("D8").Value = "0.2914"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.06772"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "113.32"
$ws.Range("E10").Value = "  +6.09%  "
$ws.Range("D11").Value = "19.32"
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").Value = "1.936.27"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "5.544"
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("D14").Value = "0.07607"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "0.6803"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "296.14"
$ws.Range("E16").Value = "  +4.12%  "
$ws.Range("D17").Value = "30.655.09"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "13.15"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "0.000007661"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.193.54"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "5.569"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "6.486"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "9.620"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "167.49"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "20.31"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "2.099"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").Value = "0.1070"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "1.427"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "4.175"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "0.04995"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").Value = "0.7536"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").Value = "1.156"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").Value = "0.02067"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").Value = "2.725"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "2.694"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "2.024"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").Value = "110.29"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Value = "0.4443"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").Value = "5.854"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "70.66"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "7.384"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "48.89"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("D48").Value = "9.314"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").Value = "0.2547"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1229"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("D51").Value = "35.19"
$ws.Range("E51").Value = "  -0.22%  "
